$wb = $excel.ActiveWorkbook

function Add-DataRow {
    param(
        [object]$ws,
        [int]$row,
        [double]$timeVal,
        [string]$bVal,
        [string]$cVal,
        [string]$dVal,
        [string]$eVal,
        [double]$fVal,
        [string]$gVal,
        [bool]$gAsText,
        [double]$hVal,
        [double]$iVal
    )

    $ws.Range("A$row").NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Range("A$row").Value = $timeVal

    $ws.Range("B$row").Value = $bVal
    $ws.Range("C$row").Value = $cVal
    $ws.Range("D$row").Value = $dVal
    $ws.Range("E$row").Value = $eVal

    $ws.Range("F$row").Value = $fVal
    if ($gAsText) {
        # Preserve the full-precision big integer as text (it would lose
        # precision if stored as a double), matching how the source data
        # represents this particular reading.
        $ws.Range("G$row").NumberFormat = "@"
    }
    $ws.Range("G$row").Value = $gVal
    $ws.Range("H$row").Value = $hVal
    $ws.Range("I$row").Value = $iVal
}

# ROW50-FE-LIFTER (sheet 1) -> append row 96
$ws1 = $wb.Worksheets.Item(1)
Add-DataRow $ws1 96 45772.30629111111 "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c," "0x01,0x3a" "0xe" 400 "568631262647113970876416" $false 314 14

# ROW50-MID-LIFTER (sheet 2) -> append row 98
$ws2 = $wb.Worksheets.Item(2)
Add-DataRow $ws2 98 45772.27034722222 "0x01,0x90 " "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c," "0x01,0x3e" "0x19" 400 "568631262647113771663628" $true 318 25

# ROW11-FE-LIFTER (sheet 3) -> append row 96
$ws3 = $wb.Worksheets.Item(3)
Add-DataRow $ws3 96 45772.33809793981 "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c," "0x01,0x3a" "0x14" 400 "568631262647113970876416" $false 314 20

# ROW11-MID-LIFTER (sheet 4) -> append row 96
$ws4 = $wb.Worksheets.Item(4)
Add-DataRow $ws4 96 45772.4573003125 "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c," "0x01,0x42" "0x19" 400 "568631262647113970876416" $false 322 25
